# Fix Training Data Issue (#48)
# The "Date" column (BF) held a malformed value ("5-27-2007-08", a mangled
# concatenation of the game date and the season) for every data row.
# Replace it with the correct ISO-style date string "2008-05-27" for each
# of the 30 team rows (rows 2-31), keeping the value as literal text
# (not an auto-converted Excel date serial).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "5-27-2007-08"
$newValue = "2008-05-27"

$firstRow = 2
$lastRow = 31
$col = "BF"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Range("$col$row")
    if ($cell.Text -eq $oldValue) {
        # Leading apostrophe forces the text "2008-05-27" to stay a plain
        # string instead of being auto-parsed into a date serial number,
        # matching the literal text that was written upstream.
        $cell.Value = "'" + $newValue
    }
}
